$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing "travel time no algo ..." columns (J,K) out to
#     make room for the two new "travel time algo3 ..." columns that were
#     tested. Read the current J/K contents first (Value2 avoids a reflection
#     quirk on the bare .Value getter in this runtime), then write them back
#     two columns to the right (L,M).
$j1 = $ws.Range("J1").Value2
$k1 = $ws.Range("K1").Value2
$j2 = $ws.Range("J2").Value2
$k2 = $ws.Range("K2").Value2
$j3 = $ws.Range("J3").Value2
$k3 = $ws.Range("K3").Value2

$ws.Range("L1").Value = $j1
$ws.Range("M1").Value = $k1
$ws.Range("L2").Value = $j2
$ws.Range("M2").Value = $k2
$ws.Range("L3").Value = $j3
$ws.Range("M3").Value = $k3

# --- Write the two new "algo3" datasets into J,K
$ws.Range("J1").Value = "travel time algo3 without HB"
$ws.Range("K1").Value = "travel time algo3 with HB"

$ws.Range("J2").Value = 1154
$ws.Range("K2").Value = 6167

$ws.Range("J3").Value = 0.74983755699999999
$ws.Range("K3").Value = 4.0071474980000001

# --- Column widths: column I loses its "best fit" flag, and the two new
#     columns (L,M) pick up the same widths the shifted columns had.
$ws.Columns.Item(9).ColumnWidth = 35.053385416666664
$ws.Columns.Item(12).ColumnWidth = 40.944010416666664
$ws.Columns.Item(13).ColumnWidth = 35.053385416666664

# --- Update the selection/view to match where the author ended up
$ws.Range("K2").Select()
$excel.ActiveWindow.ScrollColumn = 7
